# Scheduled market-data refresh: update currentAveragePrice* / Leve price & profit
# columns (H:N) across all job sheets. Values sourced from the latest Universalis
# snapshot; sheet layout/tables are unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1883.4166
$ws.Range("I38").Value = 1179.5555
$ws.Range("K38").Value = 3538.6665
$ws.Range("M38").Value = -3166.6665

$ws.Range("H129").Value = 1708.6
$ws.Range("I129").Value = 1012.8571
$ws.Range("J129").Value = 3332
$ws.Range("K129").Value = 3038.5713
$ws.Range("L129").Value = 9996
$ws.Range("M129").Value = 1961.4287
$ws.Range("N129").Value = -19996

$ws.Range("H137").Value = 3785.18
$ws.Range("I137").Value = 2549.8684
$ws.Range("J137").Value = 7697
$ws.Range("K137").Value = 7649.6052
$ws.Range("L137").Value = 23091
$ws.Range("M137").Value = -5099.6052
$ws.Range("N137").Value = -28191

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 503.44446
$ws.Range("I5").Value = 156.5
$ws.Range("J5").Value = 781
$ws.Range("K5").Value = 156.5
$ws.Range("L5").Value = 781
$ws.Range("M5").Value = -44.5
$ws.Range("N5").Value = -1005

$ws.Range("H45").Value = 22729842
$ws.Range("I45").Value = 41668356
$ws.Range("J45").Value = 3625.6
$ws.Range("K45").Value = 41668356
$ws.Range("L45").Value = 3625.6
$ws.Range("M45").Value = -41667979
$ws.Range("N45").Value = -4379.6

$ws.Range("H61").Value = 22065648
$ws.Range("I61").Value = 16135516
$ws.Range("K61").Value = 16135516
$ws.Range("M61").Value = -16135304

$ws.Range("H88").Value = 4032.7144
$ws.Range("I88").Value = 3558
$ws.Range("J88").Value = 4135.913
$ws.Range("K88").Value = 3558
$ws.Range("L88").Value = 4135.913
$ws.Range("M88").Value = -3152
$ws.Range("N88").Value = -4947.913

$ws.Range("H91").Value = 4032.7144
$ws.Range("I91").Value = 3558
$ws.Range("J91").Value = 4135.913
$ws.Range("K91").Value = 3558
$ws.Range("L91").Value = 4135.913
$ws.Range("M91").Value = -2154
$ws.Range("N91").Value = -6943.913

$ws.Range("H132").Value = 5358.2104
$ws.Range("I132").Value = 2009.3334
$ws.Range("K132").Value = 6028.0002
$ws.Range("M132").Value = -3498.0002

$ws.Range("H136").Value = 22065648
$ws.Range("I136").Value = 16135516
$ws.Range("K136").Value = 48406548
$ws.Range("M136").Value = -48403998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 503.44446
$ws.Range("I4").Value = 156.5
$ws.Range("J4").Value = 781
$ws.Range("K4").Value = 156.5
$ws.Range("L4").Value = 781
$ws.Range("M4").Value = -41.5
$ws.Range("N4").Value = -1011

$ws.Range("H94").Value = 758.26666
$ws.Range("I94").Value = 552.7778
$ws.Range("J94").Value = 1066.5
$ws.Range("K94").Value = 552.7778
$ws.Range("L94").Value = 1066.5
$ws.Range("M94").Value = -101.7778
$ws.Range("N94").Value = -1968.5

$ws.Range("H107").Value = 2690.25
$ws.Range("J107").Value = 1950
$ws.Range("L107").Value = 1950
$ws.Range("N107").Value = -5790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 838508.6
$ws.Range("I31").Value = 3177.4614
$ws.Range("J31").Value = 1562462.2
$ws.Range("K31").Value = 3177.4614
$ws.Range("L31").Value = 1562462.2
$ws.Range("M31").Value = -2882.4614
$ws.Range("N31").Value = -1563052.2

$ws.Range("H34").Value = 838508.6
$ws.Range("I34").Value = 3177.4614
$ws.Range("J34").Value = 1562462.2
$ws.Range("K34").Value = 3177.4614
$ws.Range("L34").Value = 1562462.2
$ws.Range("M34").Value = -2975.4614
$ws.Range("N34").Value = -1562866.2

$ws.Range("H58").Value = 2240.6956
$ws.Range("I58").Value = 1832.2
$ws.Range("K58").Value = 1832.2
$ws.Range("M58").Value = -1629.2

$ws.Range("H99").Value = 3750.6
$ws.Range("I99").Value = 3444.5
$ws.Range("J99").Value = 3954.6667
$ws.Range("K99").Value = 3444.5
$ws.Range("L99").Value = 3954.6667
$ws.Range("M99").Value = -1946.5
$ws.Range("N99").Value = -6950.6667

$ws.Range("H109").Value = 34249.75
$ws.Range("J109").Value = 34249.75
$ws.Range("L109").Value = 34249.75
$ws.Range("N109").Value = -36329.75

$ws.Range("H126").Value = 3750.6
$ws.Range("I126").Value = 3444.5
$ws.Range("J126").Value = 3954.6667
$ws.Range("K126").Value = 10333.5
$ws.Range("L126").Value = 11864.0001
$ws.Range("M126").Value = -7863.5
$ws.Range("N126").Value = -16804.0001

$ws.Range("H136").Value = 2240.6956
$ws.Range("I136").Value = 1832.2
$ws.Range("K136").Value = 5496.6
$ws.Range("M136").Value = -2946.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 18999.4
$ws.Range("I9").Value = 9998.5
$ws.Range("J9").Value = 25000
$ws.Range("K9").Value = 29995.5
$ws.Range("L9").Value = 75000
$ws.Range("M9").Value = -29771.5
$ws.Range("N9").Value = -75448

$ws.Range("H56").Value = 6015.4546
$ws.Range("I56").Value = 6015.4546
$ws.Range("K56").Value = 6015.4546
$ws.Range("M56").Value = -5485.4546

$ws.Range("H59").Value = 1700
$ws.Range("I59").Value = 1500
$ws.Range("J59").Value = 1800
$ws.Range("K59").Value = 4500
$ws.Range("L59").Value = 5400
$ws.Range("M59").Value = -3960
$ws.Range("N59").Value = -6480

$ws.Range("H117").Value = 1698.8
$ws.Range("J117").Value = 1698.8
$ws.Range("L117").Value = 5096.4
$ws.Range("N117").Value = -11980.4

$ws.Range("H121").Value = 2588.8572
$ws.Range("J121").Value = 3326.4
$ws.Range("L121").Value = 9979.200000000001
$ws.Range("N121").Value = -12599.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4828.5835
$ws.Range("J102").Value = 4832.1665
$ws.Range("L102").Value = 4832.1665
$ws.Range("N102").Value = -8076.1665

$ws.Range("H105").Value = 33533.4
$ws.Range("J105").Value = 33533.4
$ws.Range("L105").Value = 33533.4
$ws.Range("N105").Value = -40521.4

$ws.Range("H122").Value = 1411.4445
$ws.Range("I122").Value = 1444.9333
$ws.Range("J122").Value = 1244
$ws.Range("K122").Value = 4334.7999
$ws.Range("L122").Value = 3732
$ws.Range("M122").Value = -1884.7999
$ws.Range("N122").Value = -8632

$ws.Range("H126").Value = 3964.348
$ws.Range("I126").Value = 3356.8333
$ws.Range("K126").Value = 10070.4999
$ws.Range("M126").Value = -7600.499899999999

$ws.Range("H132").Value = 21282614
$ws.Range("I132").Value = 33336558
$ws.Range("J132").Value = 10946.412
$ws.Range("K132").Value = 100009674
$ws.Range("L132").Value = 32839.236
$ws.Range("M132").Value = -100007144
$ws.Range("N132").Value = -37899.236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

$ws.Range("H7").Value = 108570.8
$ws.Range("I7").Value = 4700
$ws.Range("K7").Value = 4700
$ws.Range("M7").Value = -4588

$ws.Range("H40").Value = 5238.769
$ws.Range("I40").Value = 4526
$ws.Range("K40").Value = 4526
$ws.Range("M40").Value = -4390

$ws.Range("H43").Value = 1262571
$ws.Range("J43").Value = 1465041.2
$ws.Range("L43").Value = 1465041.2
$ws.Range("N43").Value = -1465427.2

$ws.Range("H93").Value = 125001700
$ws.Range("J93").Value = 2042.75
$ws.Range("L93").Value = 2042.75
$ws.Range("N93").Value = -4538.75

$ws.Range("H100").Value = 4111.5454
$ws.Range("I100").Value = 3656.6
$ws.Range("K100").Value = 3656.6
$ws.Range("M100").Value = -3115.6

$ws.Range("H109").Value = 102995
$ws.Range("J109").Value = 102995
$ws.Range("L109").Value = 102995
$ws.Range("N109").Value = -105769

$ws.Range("H126").Value = 108570.8
$ws.Range("I126").Value = 4700
$ws.Range("K126").Value = 14100
$ws.Range("M126").Value = -11630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20000000
$ws.Range("J5").Value = 20000000
$ws.Range("L5").Value = 20000000
$ws.Range("N5").Value = -20000224

$ws.Range("H33").Value = 29495
$ws.Range("J33").Value = 29495
$ws.Range("L33").Value = 29495
$ws.Range("N33").Value = -29995

$ws.Range("H36").Value = 29495
$ws.Range("J36").Value = 29495
$ws.Range("L36").Value = 29495
$ws.Range("N36").Value = -29995

$ws.Range("N37").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0

$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30793

$ws.Range("H49").Value = 26747.5
$ws.Range("I49").Value = 20000
$ws.Range("K49").Value = 20000
$ws.Range("M49").Value = -19770

$ws.Range("H126").Value = 1602
$ws.Range("I126").Value = 1602
$ws.Range("K126").Value = 4806
$ws.Range("M126").Value = -2336

$ws.Range("H132").Value = 296971.53
$ws.Range("I132").Value = 2811.276
$ws.Range("K132").Value = 8433.828
$ws.Range("M132").Value = -5903.828
